# Add Doctor Details Page
# Adds new localization rows to Sheet1:
#   266: Profile / الملف الشخصي
#   267: My Profile / ملفي الشخصي   (reuses existing "My Profile" msgid)
#   268: Edit Profile / تعديل ملفي الشخصي (reuses existing "Edit Profile" msgid)
#   270: Doctor Details   (msgid only, no translation yet)
#   271: Details          (msgid only, no translation yet)
#   272: Patient Details  (msgid only, no translation yet)
#   273: Admission Details(msgid only, no translation yet)
# Row 269 is intentionally left blank/unused.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PoRow($Row, $Msgid, $Msgstr) {
    $ws.Cells.Item($Row, 1).Value = $Msgid
    $ws.Cells.Item($Row, 2).Value = $Msgstr
    $formula = '=T(_xlfn.CONCAT("msgid ", CHAR(34), ,A' + $Row + ',CHAR(34), CHAR(10), "msgstr ",CHAR(34), B' + $Row + ', CHAR(34), CHAR(10), CHAR(10)))'
    $ws.Cells.Item($Row, 3).Formula = $formula
}

# Rows with both an English source string and an Arabic translation.
Set-PoRow 266 "Profile" "الملف الشخصي"
Set-PoRow 267 "My Profile" "ملفي الشخصي"
Set-PoRow 268 "Edit Profile" "تعديل ملفي الشخصي"

# Row 269 intentionally skipped (left empty).

# Rows with only a new msgid (no translation / column B or C yet).
$ws.Cells.Item(270, 1).Value = "Doctor Details"
$ws.Cells.Item(271, 1).Value = "Details"
$ws.Cells.Item(272, 1).Value = "Patient Details"
$ws.Cells.Item(273, 1).Value = "Admission Details"

# Match the author's final selection/active cell.
$ws.Range("A273").Select() | Out-Null
